$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 03:25"

# Swap Benin / Guinea-Bisau order: Benin now comes right before Guinea-Bisau.
# Row 148 becomes Benin (with refreshed figures), row 149 becomes Guinea-Bisau
# (keeping the figures that used to belong to the old row 148 / Guinea-Bisau).
$ws.Range("A148").Value = "Benin"
$ws.Range("B148").Value = 2280
$ws.Range("C148").Value = 13
$ws.Range("D148").Value = 1942
$ws.Range("E148").Value = 298
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 40

$ws.Range("A149").Value = "Guinea-Bisau"
$ws.Range("B149").Value = 2275
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 1127
$ws.Range("E149").Value = 1109
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 39

# Refresh updated case counts for the other countries
$ws.Range("B4").Value = 6828116
$ws.Range("C4").Value = 39969
$ws.Range("D4").Value = 4118982
$ws.Range("E4").Value = 2507793
$ws.Range("G4").Value = 1144
$ws.Range("H4").Value = 201341

$ws.Range("B8").Value = 744400
$ws.Range("C8").Value = 6380
$ws.Range("D8").Value = 587717
$ws.Range("E8").Value = 125632
$ws.Range("G8").Value = 124
$ws.Range("H8").Value = 31051

$ws.Range("B29").Value = 139747
$ws.Range("C29").Value = 944
$ws.Range("D29").Value = 122449
$ws.Range("E29").Value = 8105

$ws.Range("B74").Value = 30419
$ws.Range("C74").Value = 1121
$ws.Range("D74").Value = 15740
$ws.Range("E74").Value = 14113
$ws.Range("G74").Value = 14
$ws.Range("H74").Value = 566

$ws.Range("B113").Value = 5399
$ws.Range("C113").Value = 3
$ws.Range("D113").Value = 5333
$ws.Range("E113").Value = 5

$ws.Range("B115").Value = 5155
$ws.Range("C115").Value = 27
$ws.Range("D115").Value = 4418
$ws.Range("E115").Value = 636

$ws.Range("B140").Value = 3087
$ws.Range("C140").Value = 55
$ws.Range("D140").Value = 1533
$ws.Range("E140").Value = 1485

$wb.Save()
